$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell as TEXT (matches the source inlineStr cells), avoiding
# Excel auto-converting numeric-looking strings (e.g. "5.026") into numbers,
# while leaving the cells style back at the workbook default (no explicit
# style index), matching the original unstyled data cells.
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue 'D2' '29.552.16'
$ws.Range('E2').Value = '  +1.97%  '

# Row 3
Set-TextValue 'D3' '1.842.45'
$ws.Range('E3').Value = '  +1.11%  '

# Row 4
Set-TextValue 'D4' '0.9982'
$ws.Range('E4').Value = '  -1.23%  '

# Row 5
Set-TextValue 'D5' '244.51'

# Row 6
Set-TextValue 'D6' '0.6304'
$ws.Range('E6').Value = '  +2.78%  '

# Row 7
Set-TextValue 'D7' '0.9986'
$ws.Range('E7').Value = '  -1.00%  '

# Row 8
Set-TextValue 'D8' '0.07448'
$ws.Range('E8').Value = '  +1.78%  '

# Row 9
Set-TextValue 'D9' '0.2951'
$ws.Range('E9').Value = '  +2.32%  '

# Row 10
Set-TextValue 'D10' '23.82'
$ws.Range('E10').Value = '  +4.88%  '

# Row 11
Set-TextValue 'D11' '0.07672'
$ws.Range('E11').Value = '  -0.21%  '

# Row 12
Set-TextValue 'D12' '1.841.07'
$ws.Range('E12').Value = '  +1.42%  '

# Row 13
Set-TextValue 'D13' '5.026'
$ws.Range('E13').Value = '  +2.11%  '

# Row 14
Set-TextValue 'D14' '0.6816'
$ws.Range('E14').Value = '  +3.11%  '

# Row 15
Set-TextValue 'D15' '84.13'
$ws.Range('E15').Value = '  +3.09%  '

# Row 16
Set-TextValue 'D16' '0.000009295'
$ws.Range('E16').Value = '  +3.97%  '

# Row 17
Set-TextValue 'D17' '5.975'
$ws.Range('E17').Value = '  +2.54%  '

# Row 18
Set-TextValue 'D18' '29.508.00'
$ws.Range('E18').Value = '  +1.98%  '

# Row 19
Set-TextValue 'D19' '2.081.48'
$ws.Range('E19').Value = '  +1.12%  '

# Row 20
Set-TextValue 'D20' '237.94'
$ws.Range('E20').Value = '  +1.59%  '

# Row 21
$ws.Range('E21').Value = '  +1.32%  '

# Row 22
Set-TextValue 'D22' '0.9989'
$ws.Range('E22').Value = '  -1.01%  '

# Row 23
Set-TextValue 'D23' '7.363'
$ws.Range('E23').Value = '  +4.32%  '

# Row 24
Set-TextValue 'D24' '0.9993'
$ws.Range('E24').Value = '  -1.41%  '

# Row 25
Set-TextValue 'D25' '159.03'
$ws.Range('E25').Value = '  +0.09%  '

# Row 26
Set-TextValue 'D26' '0.1418'
$ws.Range('E26').Value = '  +1.91%  '

# Row 27
Set-TextValue 'D27' '8.534'
$ws.Range('E27').Value = '  +1.44%  '

# Row 28
Set-TextValue 'D28' '17.81'
$ws.Range('E28').Value = '  +1.07%  '

# Row 29
$ws.Range('B29').Value = 'PancakeSwap'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue 'D29' '1.500'
$ws.Range('E29').Value = '  +0.66%  '

# Row 30
$ws.Range('B30').Value = 'Hedera'
$ws.Range('C30').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D30' '0.06023'
$ws.Range('E30').Value = '  +8.86%  '

# Row 31
Set-TextValue 'D31' '1.245'
$ws.Range('E31').Value = '  +3.08%  '

# Row 32
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 'D32' '4.155'
$ws.Range('E32').Value = '  +2.03%  '

# Row 33
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D33' '4.118'
$ws.Range('E33').Value = '  +0.85%  '

# Row 34
Set-TextValue 'D34' '1.867'
$ws.Range('E34').Value = '  +1.97%  '

# Row 35
Set-TextValue 'D35' '1.147'
$ws.Range('E35').Value = '  +1.62%  '

# Row 36
Set-TextValue 'D36' '0.7290'
$ws.Range('E36').Value = '  -0.20%  '

# Row 37
Set-TextValue 'D37' '2.611'
$ws.Range('E37').Value = '  -1.39%  '

# Row 38
$ws.Range('E38').Value = '  +2.30%  '

# Row 39
Set-TextValue 'D39' '1.223.67'
$ws.Range('E39').Value = '  +2.78%  '

# Row 40
$ws.Range('E40').Value = '  +0.44%  '

# Row 41
Set-TextValue 'D41' '6.280'
$ws.Range('E41').Value = '  -0.57%  '

# Row 42
Set-TextValue 'D42' '0.9174'
$ws.Range('E42').Value = '  +2.78%  '

# Row 43
Set-TextValue 'D43' '0.9999'
$ws.Range('E43').Value = '  -0.77%  '

# Row 44
Set-TextValue 'D44' '2.003.34'
$ws.Range('E44').Value = '  +2.48%  '

# Row 45
Set-TextValue 'D45' '102.06'
$ws.Range('E45').Value = '  +1.46%  '

# Row 46
Set-TextValue 'D46' '65.89'
$ws.Range('E46').Value = '  +2.97%  '

# Row 47
$ws.Range('E47').Value = '  -0.91%  '

# Row 48
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D48' '9.320'
$ws.Range('E48').Value = '  +3.99%  '

# Row 49
$ws.Range('B49').Value = 'TheSandbox'
$ws.Range('C49').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue 'D49' '0.4074'
$ws.Range('E49').Value = '  +2.32%  '

# Row 50
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue 'D50' '0.00000000118'
$ws.Range('E50').Value = '  -3.39%  '

# Row 51
Set-TextValue 'D51' '0.1138'
$ws.Range('E51').Value = '  +4.58%  '
